# UAT Items List - renumber STUDIES section (2.x -> 1.x), drop the stray
# "2" section-number cell, and add a new row for the "Study-level Consent
# Details Report" test item.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("UAT Items List")

# The "STUDIES" block used to be section "2" (A5 held the literal 2) -
# it is now section "1", so the stray numeric label in A5 is removed
# entirely (not just blanked).
$ws.Range("A5").Clear()

# Renumber the "2.x" item numbers to "1.x" (column A, rows 6-17).
$ws.Range("A6").Value = "1.1"
$ws.Range("A7").Value = "1.2"
$ws.Range("A8").Value = "1.3"
$ws.Range("A9").Value = "1.4"
$ws.Range("A10").Value = "1.5"
$ws.Range("A11").Value = "1.6"
$ws.Range("A12").Value = "1.7"
$ws.Range("A13").Value = "1.8"
$ws.Range("A14").Value = "1.9"
$ws.Range("A15").Value = "1.10"
$ws.Range("A16").Value = "1.11"
$ws.Range("A17").Value = "1.12"

# Row 19 (previously a blank follow-on row for item 1.12) now documents a
# newly discovered test item: "Study-level Consent Details Report",
# written in the same bold/green "note" style used elsewhere in the sheet.
$b19 = $ws.Range("B19")
$b19.Value = "Study-level Consent Details Report"
$b19.HorizontalAlignment = -4131
$b19.VerticalAlignment = -4160
$b19.WrapText = $true
$b19.Font.Name = "Calibri"
$b19.Font.Size = 11
$b19.Font.Bold = $true
$b19.Font.Color = 24832

# Move the on-screen selection/scroll position to reflect the newly
# added row.
$ws.Range("B19").Select() | Out-Null
$excel.ActiveWindow.ScrollRow = 13
